$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the blank separator row currently at row 12 (this shifts rows 13-18 up to 12-17)
$ws.Rows.Item(12).Delete()

# Append the new FINISH row at the bottom (row 18)
$ws.Range("A18:G18").Merge()
$ws.Range("A18").Value = "FINISH"
$ws.Range("A18:G18").HorizontalAlignment = 1

$ws.Range("D24").Select()
